$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 6
$ws_ALC.Range("H6").Value = 62.5
$ws_ALC.Range("I6").Value = 62.5
$ws_ALC.Range("K6").Value = 187.5
$ws_ALC.Range("M6").Value = -75.5

# ALC row 96
$ws_ALC.Range("H96").Value = 1159
$ws_ALC.Range("I96").Value = 1331.125
$ws_ALC.Range("J96").Value = 700
$ws_ALC.Range("K96").Value = 3993.375
$ws_ALC.Range("L96").Value = 2100
$ws_ALC.Range("M96").Value = -2620.375
$ws_ALC.Range("N96").Value = -4846

# ALC row 139
$ws_ALC.Range("H139").Value = 89874.25
$ws_ALC.Range("J139").Value = 89874.25
$ws_ALC.Range("L139").Value = 89874.25
$ws_ALC.Range("N139").Value = -100154.25

# ARM row 32
$ws_ARM.Range("H32").Value = 16047.274
$ws_ARM.Range("I32").Value = 7062.7095
$ws_ARM.Range("K32").Value = 7062.7095
$ws_ARM.Range("M32").Value = -6775.7095

# ARM row 61
$ws_ARM.Range("H61").Value = 14651.477
$ws_ARM.Range("I61").Value = 4707.357
$ws_ARM.Range("K61").Value = 4707.357
$ws_ARM.Range("M61").Value = -4495.357

# ARM row 74
$ws_ARM.Range("H74").Value = 22028
$ws_ARM.Range("I74").Value = 3177.111
$ws_ARM.Range("K74").Value = 3177.111
$ws_ARM.Range("M74").Value = -2303.111

# ARM row 77
$ws_ARM.Range("H77").Value = 22028
$ws_ARM.Range("I77").Value = 3177.111
$ws_ARM.Range("K77").Value = 15885.555
$ws_ARM.Range("M77").Value = -11517.555

# ARM row 132
$ws_ARM.Range("H132").Value = 2573878.5
$ws_ARM.Range("I132").Value = 3843.76
$ws_ARM.Range("K132").Value = 11531.28
$ws_ARM.Range("M132").Value = -9001.280000000001

# ARM row 136
$ws_ARM.Range("H136").Value = 14651.477
$ws_ARM.Range("I136").Value = 4707.357
$ws_ARM.Range("K136").Value = 14122.071
$ws_ARM.Range("M136").Value = -11572.071

# BSM row 20
$ws_BSM.Range("H20").Value = 16636.809
$ws_BSM.Range("I20").Value = 8781.913
$ws_BSM.Range("K20").Value = 8781.913
$ws_BSM.Range("M20").Value = -8534.913

# BSM row 107
$ws_BSM.Range("H107").Value = 1753.1818
$ws_BSM.Range("I107").Value = 1421.625
$ws_BSM.Range("J107").Value = 2637.3333
$ws_BSM.Range("K107").Value = 1421.625
$ws_BSM.Range("L107").Value = 2637.3333
$ws_BSM.Range("M107").Value = 498.375
$ws_BSM.Range("N107").Value = -6477.3333

# BSM row 134
$ws_BSM.Range("H134").Value = 17454.523
$ws_BSM.Range("I134").Value = 9797.727999999999
$ws_BSM.Range("K134").Value = 29393.184
$ws_BSM.Range("M134").Value = -26858.184

# CRP row 118
$ws_CRP.Range("H118").Value = 74999.336
$ws_CRP.Range("J118").Value = 74999.336
$ws_CRP.Range("L118").Value = 74999.336
$ws_CRP.Range("N118").Value = -78313.336

# CRP row 122
$ws_CRP.Range("H122").Value = 4426.517
$ws_CRP.Range("I122").Value = 2165.0527
$ws_CRP.Range("K122").Value = 6495.158100000001
$ws_CRP.Range("M122").Value = -4045.158100000001

# CUL row 26
$ws_CUL.Range("H26").Value = 1354.2222
$ws_CUL.Range("I26").Value = 1484.1428
$ws_CUL.Range("J26").Value = 899.5
$ws_CUL.Range("K26").Value = 4452.428400000001
$ws_CUL.Range("L26").Value = 2698.5
$ws_CUL.Range("M26").Value = -4164.428400000001
$ws_CUL.Range("N26").Value = -3274.5

# CUL row 44
$ws_CUL.Range("H44").Value = 417
$ws_CUL.Range("I44").Value = 125.5
$ws_CUL.Range("K44").Value = 376.5
$ws_CUL.Range("M44").Value = 21.5

# CUL row 68
$ws_CUL.Range("H68").Value = 7046.625
$ws_CUL.Range("I68").Value = 725
$ws_CUL.Range("K68").Value = 2175
$ws_CUL.Range("M68").Value = -1364

# CUL row 71
$ws_CUL.Range("H71").Value = 7046.625
$ws_CUL.Range("I71").Value = 725
$ws_CUL.Range("K71").Value = 6525
$ws_CUL.Range("M71").Value = -2469

# CUL row 87
$ws_CUL.Range("H87").Value = 16249.875
$ws_CUL.Range("I87").Value = 18333
$ws_CUL.Range("K87").Value = 54999
$ws_CUL.Range("M87").Value = -53751

# CUL row 90
$ws_CUL.Range("H90").Value = 16249.875
$ws_CUL.Range("I90").Value = 18333
$ws_CUL.Range("K90").Value = 164997
$ws_CUL.Range("M90").Value = -158757

# CUL row 129
$ws_CUL.Range("H129").Value = 2635.1428
$ws_CUL.Range("I129").Value = 3667.5
$ws_CUL.Range("J129").Value = 2222.2
$ws_CUL.Range("K129").Value = 11002.5
$ws_CUL.Range("L129").Value = 6666.599999999999
$ws_CUL.Range("M129").Value = -6002.5
$ws_CUL.Range("N129").Value = -16666.6

# CUL row 136
$ws_CUL.Range("H136").Value = 4332.6665
$ws_CUL.Range("I136").Value = 4332.6665
$ws_CUL.Range("K136").Value = 12997.9995
$ws_CUL.Range("M136").Value = -7897.999500000002

# CUL row 137
$ws_CUL.Range("H137").Value = 1496.9231
$ws_CUL.Range("I137").Value = 1076.8334
$ws_CUL.Range("K137").Value = 3230.5002
$ws_CUL.Range("M137").Value = 1869.4998

# CUL row 138
$ws_CUL.Range("H138").Value = 4349.724
$ws_CUL.Range("I138").Value = 1249.25
$ws_CUL.Range("K138").Value = 3747.75
$ws_CUL.Range("M138").Value = 1392.25

# GSM row 5
$ws_GSM.Range("H5").Value = 934.9
$ws_GSM.Range("I5").Value = 934.9
$ws_GSM.Range("K5").Value = 934.9
$ws_GSM.Range("M5").Value = -822.9

# GSM row 102
$ws_GSM.Range("H102").Value = 3452.3845
$ws_GSM.Range("I102").Value = 3822.1428
$ws_GSM.Range("J102").Value = 1899.4
$ws_GSM.Range("K102").Value = 3822.1428
$ws_GSM.Range("L102").Value = 1899.4
$ws_GSM.Range("M102").Value = -2200.1428
$ws_GSM.Range("N102").Value = -5143.4

# GSM row 122
$ws_GSM.Range("H122").Value = 2526.7778
$ws_GSM.Range("I122").Value = 2041.3334
$ws_GSM.Range("K122").Value = 6124.0002
$ws_GSM.Range("M122").Value = -3674.0002

# GSM row 132
$ws_GSM.Range("H132").Value = 10818.137
$ws_GSM.Range("I132").Value = 6929.2144
$ws_GSM.Range("K132").Value = 20787.6432
$ws_GSM.Range("M132").Value = -18257.6432

# GSM row 136
$ws_GSM.Range("H136").Value = 7499.5
$ws_GSM.Range("J136").Value = 7499.5
$ws_GSM.Range("L136").Value = 22498.5
$ws_GSM.Range("N136").Value = -27598.5

# GSM row 140
$ws_GSM.Range("H140").Value = 67855.7
$ws_GSM.Range("J140").Value = 67855.7
$ws_GSM.Range("L140").Value = 67855.7
$ws_GSM.Range("N140").Value = -78215.7

# LTW row 7
$ws_LTW.Range("H7").Value = 9999.733
$ws_LTW.Range("I7").Value = 4666.1665
$ws_LTW.Range("K7").Value = 4666.1665
$ws_LTW.Range("M7").Value = -4554.1665

# LTW row 16
$ws_LTW.Range("H16").Value = 2666.7273
$ws_LTW.Range("I16").Value = 2728.3333
$ws_LTW.Range("J16").Value = 2389.5
$ws_LTW.Range("K16").Value = 2728.3333
$ws_LTW.Range("L16").Value = 2389.5
$ws_LTW.Range("M16").Value = -2558.3333
$ws_LTW.Range("N16").Value = -2729.5

# LTW row 61
$ws_LTW.Range("H61").Value = 2492.8
$ws_LTW.Range("I61").Value = 1604.8667
$ws_LTW.Range("K61").Value = 1604.8667
$ws_LTW.Range("M61").Value = -1402.8667

# LTW row 93
$ws_LTW.Range("H93").Value = 18074.166
$ws_LTW.Range("I93").Value = 14127.571
$ws_LTW.Range("J93").Value = 23599.4
$ws_LTW.Range("K93").Value = 14127.571
$ws_LTW.Range("L93").Value = 23599.4
$ws_LTW.Range("M93").Value = -12879.571
$ws_LTW.Range("N93").Value = -26095.4

# LTW row 113
$ws_LTW.Range("H113").Value = 2492.8
$ws_LTW.Range("I113").Value = 1604.8667
$ws_LTW.Range("K113").Value = 1604.8667
$ws_LTW.Range("M113").Value = 565.1333

# LTW row 126
$ws_LTW.Range("H126").Value = 9999.733
$ws_LTW.Range("I126").Value = 4666.1665
$ws_LTW.Range("K126").Value = 13998.4995
$ws_LTW.Range("M126").Value = -11528.4995

# LTW row 132
$ws_LTW.Range("H132").Value = 4012687
$ws_LTW.Range("I132").Value = 3358.25
$ws_LTW.Range("K132").Value = 10074.75
$ws_LTW.Range("M132").Value = -7544.75

# WVR row 62
$ws_WVR.Range("H62").Value = 2999.8
$ws_WVR.Range("J62").Value = 2999.6667
$ws_WVR.Range("L62").Value = 2999.6667
$ws_WVR.Range("N62").Value = -4247.6667

# WVR row 65
$ws_WVR.Range("H65").Value = 2999.8
$ws_WVR.Range("J65").Value = 2999.6667
$ws_WVR.Range("L65").Value = 14998.3335
$ws_WVR.Range("N65").Value = -21238.3335

# WVR row 126
$ws_WVR.Range("H126").Value = 25605.938
$ws_WVR.Range("I126").Value = 23979.732
$ws_WVR.Range("K126").Value = 71939.196
$ws_WVR.Range("M126").Value = -69469.196

# WVR row 132
$ws_WVR.Range("H132").Value = 16869.7
$ws_WVR.Range("I132").Value = 1527.5714
$ws_WVR.Range("K132").Value = 4582.7142
$ws_WVR.Range("M132").Value = -2052.7142

# WVR row 136
$ws_WVR.Range("H136").Value = 14695.526
$ws_WVR.Range("I136").Value = 3243.7144
$ws_WVR.Range("K136").Value = 9731.143199999999
$ws_WVR.Range("M136").Value = -7181.143199999999

# WVR row 141
$ws_WVR.Range("H141").Value = 75749.75
$ws_WVR.Range("J141").Value = 75749.75
$ws_WVR.Range("L141").Value = 75749.75
$ws_WVR.Range("N141").Value = -86109.75
